$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 149.16667
$ws.Range("I33").Value = 163.42857
$ws.Range("K33").Value = 163.42857
$ws.Range("M33").Value = 65.57142999999999
$ws.Range("H34").Value = 2435.5
$ws.Range("I34").Value = 2435.5
$ws.Range("K34").Value = 2435.5
$ws.Range("M34").Value = -2232.5
$ws.Range("H36").Value = 2435.5
$ws.Range("I36").Value = 2435.5
$ws.Range("K36").Value = 2435.5
$ws.Range("M36").Value = -1720.5
$ws.Range("H106").Value = 33336578
$ws.Range("I106").Value = 55556972
$ws.Range("K106").Value = 55556972
$ws.Range("M106").Value = -55556341
$ws.Range("H107").Value = 656.6667
$ws.Range("I107").Value = 656.6667
$ws.Range("K107").Value = 656.6667
$ws.Range("M107").Value = 1263.3333
$ws.Range("H132").Value = 6528.34
$ws.Range("I132").Value = 1413.1228
$ws.Range("J132").Value = 13308.977
$ws.Range("K132").Value = 4239.3684
$ws.Range("L132").Value = 39926.931
$ws.Range("M132").Value = -1709.3684
$ws.Range("N132").Value = -44986.931
$ws.Range("H137").Value = 5619.231
$ws.Range("I137").Value = 2499.75
$ws.Range("K137").Value = 7499.25
$ws.Range("M137").Value = -4949.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5983050.5
$ws.Range("I2").Value = 8724576
$ws.Range("K2").Value = 8724576
$ws.Range("M2").Value = -8724463
$ws.Range("H4").Value = 350
$ws.Range("I4").Value = 350
$ws.Range("K4").Value = 350
$ws.Range("M4").Value = -234
$ws.Range("H5").Value = 738.25
$ws.Range("I5").Value = 738.25
$ws.Range("K5").Value = 738.25
$ws.Range("M5").Value = -626.25
$ws.Range("H32").Value = 3525.0833
$ws.Range("I32").Value = 1574.8695
$ws.Range("K32").Value = 1574.8695
$ws.Range("M32").Value = -1287.8695
$ws.Range("H74").Value = 1701.5714
$ws.Range("I74").Value = 1303.6666
$ws.Range("K74").Value = 1303.6666
$ws.Range("M74").Value = -429.6666
$ws.Range("H77").Value = 1701.5714
$ws.Range("I77").Value = 1303.6666
$ws.Range("K77").Value = 6518.333000000001
$ws.Range("M77").Value = -2150.333000000001
$ws.Range("H116").Value = 5983050.5
$ws.Range("I116").Value = 8724576
$ws.Range("K116").Value = 8724576
$ws.Range("M116").Value = -8722282
$ws.Range("H122").Value = 5429.5713
$ws.Range("I122").Value = 1335.6666
$ws.Range("K122").Value = 4006.9998
$ws.Range("M122").Value = -1556.9998
$ws.Range("H132").Value = 11356.141
$ws.Range("J132").Value = 6398.278
$ws.Range("L132").Value = 19194.834
$ws.Range("N132").Value = -24254.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5983050.5
$ws.Range("I3").Value = 8724576
$ws.Range("K3").Value = 8724576
$ws.Range("M3").Value = -8724462
$ws.Range("H4").Value = 738.25
$ws.Range("I4").Value = 738.25
$ws.Range("K4").Value = 738.25
$ws.Range("M4").Value = -623.25
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H37").Value = 962.75
$ws.Range("I37").Value = 950.6667
$ws.Range("J37").Value = 999
$ws.Range("K37").Value = 950.6667
$ws.Range("L37").Value = 999
$ws.Range("M37").Value = -813.6667
$ws.Range("N37").Value = -1273
$ws.Range("H86").Value = 2755.2
$ws.Range("I86").Value = 2465.353
$ws.Range("K86").Value = 2465.353
$ws.Range("M86").Value = -1342.353
$ws.Range("H89").Value = 2755.2
$ws.Range("I89").Value = 2465.353
$ws.Range("K89").Value = 12326.765
$ws.Range("M89").Value = -6710.764999999999
$ws.Range("H102").Value = 39995.6
$ws.Range("J102").Value = 44994.75
$ws.Range("L102").Value = 44994.75
$ws.Range("N102").Value = -51484.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3728.0876
$ws.Range("I31").Value = 3040.468
$ws.Range("J31").Value = 6959.9
$ws.Range("K31").Value = 3040.468
$ws.Range("L31").Value = 6959.9
$ws.Range("M31").Value = -2745.468
$ws.Range("N31").Value = -7549.9
$ws.Range("H34").Value = 3728.0876
$ws.Range("I34").Value = 3040.468
$ws.Range("J34").Value = 6959.9
$ws.Range("K34").Value = 3040.468
$ws.Range("L34").Value = 6959.9
$ws.Range("M34").Value = -2838.468
$ws.Range("N34").Value = -7363.9
$ws.Range("H58").Value = 401358.38
$ws.Range("I58").Value = 556752.6
$ws.Range("K58").Value = 556752.6
$ws.Range("M58").Value = -556549.6
$ws.Range("H134").Value = 1860.86
$ws.Range("I134").Value = 1617.7954
$ws.Range("J134").Value = 3643.3333
$ws.Range("K134").Value = 4853.3862
$ws.Range("L134").Value = 10929.9999
$ws.Range("M134").Value = -2318.3862
$ws.Range("N134").Value = -15999.9999
$ws.Range("H136").Value = 401358.38
$ws.Range("I136").Value = 556752.6
$ws.Range("K136").Value = 1670257.8
$ws.Range("M136").Value = -1667707.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 1182129.4
$ws.Range("K4").Value = 3546388.2
$ws.Range("M4").Value = -3546276.2
$ws.Range("H6").Value = 320
$ws.Range("I6").Value = 320
$ws.Range("K6").Value = 960
$ws.Range("M6").Value = -847
$ws.Range("H9").Value = 709
$ws.Range("I9").Value = 820
$ws.Range("K9").Value = 2460
$ws.Range("M9").Value = -2236
$ws.Range("H12").Value = 203.86957
$ws.Range("J12").Value = 179.41176
$ws.Range("L12").Value = 538.23528
$ws.Range("N12").Value = -884.23528
$ws.Range("H41").Value = 500
$ws.Range("I41").Value = 500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -1162
$ws.Range("N41").ClearContents()
$ws.Range("H122").Value = 467.93332
$ws.Range("J122").Value = 467.58334
$ws.Range("L122").Value = 4208.25006
$ws.Range("N122").Value = -9108.25006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H102").Value = 5012.9536
$ws.Range("I102").Value = 4532.1943
$ws.Range("J102").Value = 7485.4287
$ws.Range("K102").Value = 4532.1943
$ws.Range("L102").Value = 7485.4287
$ws.Range("M102").Value = -2910.1943
$ws.Range("N102").Value = -10729.4287

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5632.9
$ws.Range("I132").Value = 4499.6
$ws.Range("J132").Value = 7332.85
$ws.Range("K132").Value = 13498.8
$ws.Range("L132").Value = 21998.55
$ws.Range("M132").Value = -10968.8
$ws.Range("N132").Value = -27058.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1215.0588
$ws.Range("I126").Value = 1255.8889
$ws.Range("J126").Value = 1169.125
$ws.Range("K126").Value = 3767.6667
$ws.Range("L126").Value = 3507.375
$ws.Range("M126").Value = -1297.6667
$ws.Range("N126").Value = -8447.375
$ws.Range("H132").Value = 24162594
$ws.Range("I132").Value = 3970113
$ws.Range("K132").Value = 11910339
$ws.Range("M132").Value = -11907809
